$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F3 previously held "EAP - MCT-2A" - it moves to F6, so clear F3
$ws.Range("F3").Value = "-"

# Row 6 previously had "EAP - MEC-3A" in C6; it moves to E6, and
# "EAP - MCT-2A" (from F3) moves into F6. C6 and D6 become "-".
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "EAP - MEC-3A"
$ws.Range("F6").Value = "EAP - MCT-2A"
